$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "56.378.92"
$ws.Range("E2").Value = "  +3.97%  "

# Row 3
$ws.Range("D3").Value = "2.316.58"
$ws.Range("E3").Value = "  +2.31%  "

# Row 4
$ws.Range("D4").Value = "'1.01"
$ws.Range("E4").Value = "  +0.57%  "

# Row 5
$ws.Range("D5").Value = "'518.37"
$ws.Range("E5").Value = "  +4.49%  "

# Row 6
$ws.Range("D6").Value = "'134.08"
$ws.Range("E6").Value = "  +4.08%  "

# Row 7
$ws.Range("D7").Value = "'0.997"
$ws.Range("E7").Value = "  -0.10%  "

# Row 8
$ws.Range("D8").Value = "'0.537"
$ws.Range("E8").Value = "  +2.25%  "

# Row 9
$ws.Range("D9").Value = "2.342.56"
$ws.Range("E9").Value = "  +3.25%  "

# Row 10
$ws.Range("E10").Value = "  +9.00%  "

# Row 11
$ws.Range("D11").Value = "'0.154"
$ws.Range("E11").Value = "  +1.06%  "

# Row 12
$ws.Range("D12").Value = "'5.16"
$ws.Range("E12").Value = "  +6.99%  "

# Row 13
$ws.Range("E13").Value = "  +2.76%  "

# Row 14
$ws.Range("D14").Value = "'24.17"
$ws.Range("E14").Value = "  +5.49%  "

# Row 15
$ws.Range("D15").Value = "2.730.65"
$ws.Range("E15").Value = "  +2.52%  "

# Row 16
$ws.Range("D16").Value = "56.605.69"
$ws.Range("E16").Value = "  +4.47%  "

# Row 17
$ws.Range("E17").Value = "  +4.85%  "

# Row 18
$ws.Range("D18").Value = "2.332.73"
$ws.Range("E18").Value = "  +2.88%  "

# Row 19
$ws.Range("E19").Value = "  +3.49%  "

# Row 20
$ws.Range("D20").Value = "'4.29"
$ws.Range("E20").Value = "  +4.12%  "

# Row 21
$ws.Range("D21").Value = "'321.81"
$ws.Range("E21").Value = "  +6.20%  "

# Row 22
$ws.Range("D22").Value = "'6.65"
$ws.Range("E22").Value = "  +5.05%  "

# Row 23
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  -0.10%  "

# Row 24
$ws.Range("D24").Value = "'61.29"
$ws.Range("E24").Value = "  +1.02%  "

# Row 25
$ws.Range("D25").Value = "'0.992"
$ws.Range("E25").Value = "  -0.84%  "

# Row 26
$ws.Range("D26").Value = "'0.159"
$ws.Range("E26").Value = "  +6.83%  "

# Row 27
$ws.Range("D27").Value = "'7.69"
$ws.Range("E27").Value = "  +5.65%  "

# Row 28
$ws.Range("D28").Value = "'171.68"
$ws.Range("E28").Value = "  +0.13%  "

# Row 29
$ws.Range("D29").Value = "'1.21"
$ws.Range("E29").Value = "  +12.43%  "

# Row 30
$ws.Range("D30").Value = "0.0₃0736"
$ws.Range("E30").Value = "  +6.88%  "

# Row 31
$ws.Range("B31").Value = "Aptos"
$ws.Range("C31").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D31").Value = "'6.28"
$ws.Range("E31").Value = "  +5.28%  "

# Row 32
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "'1.68"
$ws.Range("E32").Value = "  +4.87%  "

# Row 33
$ws.Range("D33").Value = "'18.42"
$ws.Range("E33").Value = "  +3.68%  "

# Row 34
$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "  -0.02%  "

# Row 35
$ws.Range("B35").Value = "SuiNetwork"
$ws.Range("C35").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D35").Value = "'0.958"
$ws.Range("E35").Value = "  +1.49%  "

# Row 36
$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D36").Value = "'0.992"
$ws.Range("E36").Value = "  -0.45%  "

# Row 37
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'1.27"
$ws.Range("E37").Value = "  +6.09%  "

# Row 38
$ws.Range("D38").Value = "'4.04"
$ws.Range("E38").Value = "  +9.23%  "

# Row 39
$ws.Range("D39").Value = "'1.52"
$ws.Range("E39").Value = "  +9.02%  "

# Row 40
$ws.Range("D40").Value = "'37.57"
$ws.Range("E40").Value = "  +4.69%  "

# Row 41
$ws.Range("E41").Value = "  +2.30%  "

# Row 42
$ws.Range("D42").Value = "'140.29"
$ws.Range("E42").Value = "  +12.44%  "

# Row 43
$ws.Range("E43").Value = "  +7.14%  "

# Row 44
$ws.Range("D44").Value = "'5.15"
$ws.Range("E44").Value = "  +7.26%  "

# Row 45
$ws.Range("D45").Value = "'276.14"
$ws.Range("E45").Value = "  +14.48%  "

# Row 46
$ws.Range("D46").Value = "'0.0511"
$ws.Range("E46").Value = "  +3.82%  "

# Row 48
$ws.Range("D48").Value = "'0.563"
$ws.Range("E48").Value = "  +3.41%  "

# Row 49
$ws.Range("E49").Value = "  +6.41%  "

# Row 50
$ws.Range("E50").Value = "  +2.29%  "

# Row 51
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'16.95"
$ws.Range("E51").Value = "  +5.35%  "
